# Update the generated-sentence similarity metrics on each study case /
# experiment sheet, and correct one language-pair label on "Study case 2".

$wb = $excel.ActiveWorkbook

# Sheet "Study case 1" (sheet1.xml)
$ws = $wb.Worksheets.Item("Study case 1")
$ws.Range("B2").Value = 0.7592708010895827
$ws.Range("C2").Value = 4.143796545088468
$ws.Range("B3").Value = 0.7471509977842862
$ws.Range("C3").Value = 4.229964492103247

# Sheet "Study case 2" (sheet2.xml) - also relabels rows 2 and 3
$ws = $wb.Worksheets.Item("Study case 2")
$ws.Range("A2").Value = "en -> de"
$ws.Range("B2").Value = 0.749783041539854
$ws.Range("C2").Value = 4.820215035038566
$ws.Range("A3").Value = "en -> sv -> de"
$ws.Range("B3").Value = 0.737107234633723
$ws.Range("C3").Value = 4.912853220126173

# Sheet "Study case 3" (sheet3.xml)
$ws = $wb.Worksheets.Item("Study case 3")
$ws.Range("B2").Value = 0.7592708010895827
$ws.Range("C2").Value = 4.143796545088468
$ws.Range("B3").Value = 0.7323715599868168
$ws.Range("C3").Value = 4.375050645321163

# Sheet "Study case 4" (sheet4.xml)
$ws = $wb.Worksheets.Item("Study case 4")
$ws.Range("B2").Value = 0.7592708010895827
$ws.Range("C2").Value = 4.143796545088468
$ws.Range("B3").Value = 0.7323715599868168
$ws.Range("C3").Value = 4.375050645321163

# Sheet "Experiment 1" (sheet5.xml)
$ws = $wb.Worksheets.Item("Experiment 1")
$ws.Range("B2").Value = 0.668403364983767
$ws.Range("C2").Value = 4.698783337542025
$ws.Range("B3").Value = 0.6620158784065271
$ws.Range("C3").Value = 4.732067402108964

# Sheet "Experiment 2" (sheet6.xml)
$ws = $wb.Worksheets.Item("Experiment 2")
$ws.Range("B2").Value = 0.668403364983767
$ws.Range("C2").Value = 4.698783337542025
$ws.Range("B3").Value = 0.6537497597011482
$ws.Range("C3").Value = 4.785809072147154

# Sheet "Experiment 3" (sheet7.xml)
$ws = $wb.Worksheets.Item("Experiment 3")
$ws.Range("B2").Value = 0.7066478643152803
$ws.Range("C2").Value = 5.405672186625698
$ws.Range("B3").Value = 0.6900929110843625
$ws.Range("C3").Value = 5.535204230993649

# Sheet "Experiment 4" (sheet8.xml)
$ws = $wb.Worksheets.Item("Experiment 4")
$ws.Range("B2").Value = 0.7066478643152803
$ws.Range("C2").Value = 5.405672186625698
$ws.Range("B3").Value = 0.6993193104625874
$ws.Range("C3").Value = 5.459112877717812
